$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.062288
$ws.Range("H2").Value = 0.186864
$ws.Range("I2").Value = 0.01284041117846354
$ws.Range("J2").Value = 0.01284041117846353
$ws.Range("Q2").Value = 0.05206155616
$ws.Range("R2").Value = 0.46855400544
$ws.Range("S2").Value = 0.01284041117846354
$ws.Range("T2").Value = 0.01284041117846353

# Row 3
$ws.Range("I3").Value = 0.07774638634957316
$ws.Range("J3").Value = 0.07774638634957316
$ws.Range("S3").Value = 0.07774638634957316
$ws.Range("T3").Value = 0.07774638634957316

# Row 4
$ws.Range("G4").Value = 4.411504333333333
$ws.Range("I4").Value = 0.9094132024719633
$ws.Range("J4").Value = 0.9094132024719633
$ws.Range("R4").Value = 33.18501196698
$ws.Range("S4").Value = 0.9094132024719633
$ws.Range("T4").Value = 0.9094132024719633
